$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Sheet1"

$values = @(
    'recovery compression shirt',
    'recovery compression shorts',
    'recovery compression tights',
    'recovery compression tights men',
    'recovery compression tights women',
    'recovery compression top men',
    'recovery games for addiction',
    'recovery leggings',
    'recovery leggings men',
    'recovery pants',
    'recovery pants men',
    'recovery pants women',
    'recovery ring',
    'recovery shakes for women',
    'recovery shorts',
    'recovery tape',
    'recovery tights women',
    'red and black yoga pants for women',
    'red and green tights for women',
    'red and grey tights',
    'red basketball compression pants',
    'red basketball leggings',
    'red black tights women',
    'red boxing tape',
    'red capri leggings plus',
    'red capris leggings',
    'red cigarette pants for women',
    'red climbing pants',
    'red compression pants',
    'red compression tights women',
    'red football leggings for men',
    'red football tights for men',
    'red leggings low rise',
    'red leggings men',
    'red leggins for womens',
    'red long leggings',
    'red lululemon leggings',
    'red lycra pants',
    'red lycra shorts',
    'red nylon pantyhose',
    'red pants women high waist',
    'red running capris',
    'red running leggings womens',
    'red running pants women',
    'red running skirts for women',
    'red sheer tights',
    'red skins football gear',
    'red snowboarding pants women',
    'red tape sports',
    'red team vs blue team',
    'red thigh tights',
    'red tights for women',
    'red tights with control top',
    'red tights women',
    'red vs blue game',
    'red weave on track',
    'red white and blue riding pants',
    'red white blue leggings',
    'red workout legging',
    'red zebra top',
    'reduce muscle fatigue',
    'reebok activewear',
    'reebok clothes women',
    'reebok clothing women',
    'reebok compression tights',
    'reebok crossfit',
    'reebok crossfit compression',
    'reebok crossfit women',
    'reebok fitness',
    'reebok joggers women',
    'reebok lux high rise tight',
    'reebok mid',
    'reebok tights women',
    'reebok trail',
    'reebok women',
    'reebok women clothing',
    'reebok women crossfit',
    'reebok women running tights',
    'reebok womens clothing',
    'reebok womens compression tights',
    'reebok workout plus',
    'reefer parts',
    'reflective apparel',
    'reflective ball',
    'reflective bike clothing',
    'reflective capri',
    'reflective clothes',
    'reflective clothing for runners',
    'reflective clothing for women',
    'reflective clothing running',
    'reflective clothing women',
    'reflective cold weather gear',
    'reflective designs',
    'reflective fabric tape for clothing',
    'reflective gear biking',
    'reflective gear for biking',
    'reflective gear for runners',
    'reflective gear women',
    'reflective hoop',
    'reflective jacket for women',
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
